$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 100
$ws.Range("B2").Value = 106
$ws.Range("C2").ClearContents()

$ws.Range("A4").Value = 160
$ws.Range("B4").Value = 167
$ws.Range("C4").ClearContents()

$ws.Range("B6").Value = -400
$ws.Range("C6").ClearContents()

$ws.Range("B7").Value = 430
$ws.Range("C7").ClearContents()

$ws.Range("P22").Formula = "=118*1.14"
$ws.Range("P23").Formula = "=P22*15"

$ws.Range("B5").Select()
